$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (14th column), shifting existing
# columns N..S to O..T. This makes room for the new "Kindergarten" flag
# between "Vorschulkind" and "Schulkind".
$ws.Columns.Item(14).Insert()

# Match the width of the column immediately to the left (column M), as Excel
# does when inserting a column in the middle of a formatted range.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Populate the newly inserted column N with the header/value placeholders.
$ws.Range("N4").Value = "{kindergartenTitle}"
$ws.Range("N5").Value = "{isKindergarten}"
